$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("NG", "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/gpu-accelerated/ndmi300xv5-series?tabs=sizebasic"),
    @("NV", "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/gpu-accelerated/nvv3-series?tabs=sizebasic"),
    @($null, "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/gpu-accelerated/nvv4-series?tabs=sizebasic"),
    @($null, "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/gpu-accelerated/nvadsa10v5-series?tabs=sizebasic"),
    @($null, "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/gpu-accelerated/nvadsv710-v5-series?tabs=sizebasic"),
    @("NM", "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/fpga-accelerated/np-series?tabs=sizebasic"),
    @($null, "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/fpga-accelerated/nm-ads-ma35d-series?tabs=Basics"),
    @("HB", "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/high-performance-compute/hb-series?tabs=sizebasic"),
    @($null, "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/high-performance-compute/hbv2-series?tabs=sizebasic"),
    @($null, "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/high-performance-compute/hbv3-series?tabs=sizebasic"),
    @($null, "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/high-performance-compute/hbv4-series?tabs=sizebasic"),
    @("HC", "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/high-performance-compute/hc-family"),
    @($null, "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/high-performance-compute/hc-series?tabs=sizebasic"),
    @($null, "https://learn.microsoft.com/de-de/azure/virtual-machines/sizes/high-performance-compute/hx-series?tabs=sizebasic")
)

$startRow = 93
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $a = $data[$i][0]
    $b = $data[$i][1]
    if ($a -ne $null) {
        $ws.Cells.Item($row, 1).Value = $a
    }
    $ws.Cells.Item($row, 2).Value = $b
}

$ws.Range("B107").Select()
